$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-19"

# Update the header label for the Total column
$ws.Range("I1").Value = "2022 (through 07-19)"

# Update the data values for the new day's data (2022-07-27 commit, data through 07-19)
$ws.Range("I7").Value = 143
$ws.Range("I8").Value = 105
$ws.Range("I14").Value = 911
